$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) --------------------------------------------------
# Rename existing headers (G1:K1) in place.
$ws.Range("G1").Value = "점수(룰)"
$ws.Range("H1").Value = "3일상승확률(%)"
$ws.Range("I1").Value = "5일상승확률(%)"
$ws.Range("J1").Value = "10일상승확률(%)"
$ws.Range("K1").Value = "최종점수"

# New header cells L1:O1 - copy formatting from the existing header style
# (bold / bordered / centered) before writing the new text.
$ws.Range("K1").Copy()
$ws.Range("L1:O1").PasteSpecial(-4122)
$ws.Range("L1").Value = "예측방식"
$ws.Range("M1").Value = "판단"
$ws.Range("N1").Value = "MACRO_SCORE"
$ws.Range("O1").Value = "MACRO_SIGNAL"

# --- Row 2 (now Archer Aviation / ACHR) -----------------------------------
# Force the date-looking string to stay plain text (as in the source file)
# instead of Excel auto-converting it to a date serial number.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2025-11-29"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = "Archer Aviation Inc."
$ws.Range("C2").Value = "ACHR"
$ws.Range("D2").Value = 7.82
$ws.Range("E2").ClearContents()
$ws.Range("F2").Value = 13.76
$ws.Range("I2").Value = 63
$ws.Range("K2").Value = 59.8
$ws.Range("L2").Value = "Pattern"
$ws.Range("M2").Value = "⛔ 관망하십시오."
$ws.Range("N2").Value = 85.36763896678245
$ws.Range("O2").Value = "🟢 완화적 (상승 우위)"

# --- Row 3 (now Joby Aviation / JOBY) -------------------------------------
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2025-11-29"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = "Joby Aviation, Inc."
$ws.Range("C3").Value = "JOBY"
$ws.Range("D3").Value = 14.43
$ws.Range("E3").ClearContents()
$ws.Range("F3").Value = 10.62
$ws.Range("I3").Value = 63
$ws.Range("J3").Value = 66
$ws.Range("K3").Value = 56.8
$ws.Range("L3").Value = "Pattern"
$ws.Range("M3").Value = "⛔ 관망하십시오."
$ws.Range("N3").Value = 85.36763896678245
$ws.Range("O3").Value = "🟢 완화적 (상승 우위)"
